$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 216; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value -eq 46061) {
        $cell.Value = 46062
    }
}
